# Add "Data Arquivo Ori" column (F) with a fixed timestamp value for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 6).Value = "Data Arquivo Ori"

for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 6).Value = "2024-11-28 09:45:01"
}
